# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" message text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.34 = 8856.23 pesos`n✅ 8856.23 pesos = 2.33 = 948.53 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Hoja "tasas": update Binance rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 427
$ws2.Range("O10").Value = 3781.61
$ws2.Range("N12").Value = 3801
$ws2.Range("O12").Value = 407.1
